# Case study 3 inst
# Adds two new slides (using the "Title and Content" layout, i.e. ppLayoutText)
# after the existing title slide:
#   Slide 2: "Clustering notes"
#   Slide 3: "Clustering questions "

$p = $ppt.ActivePresentation

# --- Slide 2: "Clustering notes" ---
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Title.TextFrame.TextRange.Text = "Clustering notes"

# --- Slide 3: "Clustering questions " ---
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Title.TextFrame.TextRange.Text = "Clustering questions "
